$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2021 column (R) mirroring the existing 2020 column (Q) formatting
$ws.Range("R2").Style = $ws.Range("Q2").Style
$ws.Range("R3").Style = $ws.Range("Q3").Style
$ws.Range("R4").Style = $ws.Range("Q4").Style
$ws.Range("R5").Style = $ws.Range("Q5").Style

$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 13.5
$ws.Range("R5").Value = 15.1

# Update selection to match target state
$ws.Range("T3").Select()
